$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.714.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.803.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5948"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2783"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06857"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07553"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.809.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.731"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6283"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.049.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009306"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.703.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.486"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.872"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.856"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1275"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.451"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06255"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.421"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.779"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.764"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.720"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.058"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.502"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.726"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.426"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.142.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8670"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.968.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.67%  "

$ws.Range("E47").Value = "  -5.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.593"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.352"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05463"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4496"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.38%  "

